# Add column for "Status as of July 4, 2025" with an in-cell dropdown
# sourced from a new hidden helper sheet "DropdownOptions".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. Create the hidden helper sheet right after Sheet1 -----------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "DropdownOptions"

$options = @("0% - 10%", "11% - 25%", "26% - 50%", "51% - 75%", "76% - 90%", "91% - 99%", "100%")
for ($i = 0; $i -lt $options.Length; $i++) {
    $cell = $ws2.Cells.Item($i + 1, 1)
    # Force text so values like "100%" aren't coerced into a percentage number,
    # then drop the explicit style so the cell stays format-less.
    $cell.NumberFormat = "@"
    $cell.Value = $options[$i]
    $cell.ClearFormats()
}

$ws2.Visible = $false

# --- 2. Add the new header column on Sheet1 --------------------------------
$ws1.Range("AU1").Value = "Status as of July 4, 2025"

# --- 3. Apply the dropdown list validation to AU2:AU15 ---------------------
$validation = $ws1.Range("AU2:AU15").Validation
$validation.Add(3, 1, 1, '=DropdownOptions!$A$1:$A$7')
$validation.IgnoreBlank = $true
$validation.InCellDropdown = $true
$validation.ShowInput = $false
$validation.ShowError = $false

# --- 4. Leave the workbook focused back on the main sheet -----------------
$ws1.Activate()
$ws1.Range("A1").Select() | Out-Null
